# Modificando horas del mes
# Update the hours worked on row 11 of the timesheet (Sheet1) so that every
# worked day shows 8 hours instead of the previous varying values.
# Downstream SUM formulas (row 18 daily totals, AE23 weekly total, AE30
# monthly total) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newHours = @{
    "B11" = 8
    "C11" = 8
    "D11" = 8
    "E11" = 8
    "H11" = 8
    "I11" = 8
    "J11" = 8
    "K11" = 8
    "L11" = 8
    "O11" = 8
    "P11" = 8
    "Q11" = 8
    "S11" = 8
    "V11" = 8
    "W11" = 8
    "X11" = 8
    "Z11" = 8
    "AC11" = 8
    "AD11" = 8
}

foreach ($addr in $newHours.Keys) {
    $ws.Range($addr).Value = $newHours[$addr]
}

$excel.Calculate()

# Reflect the cell that ended up selected after the edit.
$ws.Range("AG11").Select()
